$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 368.7
$ws.Range("I98").Value = 327.94116
$ws.Range("J98").Value = 599.6667
$ws.Range("K98").Value = 327.94116
$ws.Range("L98").Value = 599.6667
$ws.Range("M98").Value = 1170.05884
$ws.Range("N98").Value = -3595.6667
$ws.Range("H101").Value = 11666
$ws.Range("I101").Value = 13000.5
$ws.Range("J101").Value = 990
$ws.Range("K101").Value = 39001.5
$ws.Range("L101").Value = 2970
$ws.Range("M101").Value = -37379.5
$ws.Range("N101").Value = -6214
$ws.Range("H121").Value = 1900
$ws.Range("I121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("M121").Value = ""
$ws.Range("H122").Value = 368.7
$ws.Range("I122").Value = 327.94116
$ws.Range("J122").Value = 599.6667
$ws.Range("K122").Value = 983.82348
$ws.Range("L122").Value = 1799.0001
$ws.Range("M122").Value = 1466.17652
$ws.Range("N122").Value = -6699.0001
$ws.Range("H141").Value = 2764.2307
$ws.Range("I141").Value = 2058.889
$ws.Range("J141").Value = 4351.25
$ws.Range("K141").Value = 6176.667
$ws.Range("L141").Value = 13053.75
$ws.Range("M141").Value = -996.6670000000004
$ws.Range("N141").Value = -23413.75

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 53666.668
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 53666.668
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 53666.668
$ws.Range("M139").Value = ""
$ws.Range("N139").Value = -63946.668

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 35900
$ws.Range("J140").Value = 35900
$ws.Range("L140").Value = 35900
$ws.Range("N140").Value = -46260

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 7919
$ws.Range("J51").Value = 7900
$ws.Range("L51").Value = 7900
$ws.Range("N51").Value = -9372
$ws.Range("H61").Value = 7919
$ws.Range("J61").Value = 7900
$ws.Range("L61").Value = 7900
$ws.Range("N61").Value = -8596
$ws.Range("H68").Value = 14422.385
$ws.Range("J68").Value = 14422.385
$ws.Range("L68").Value = 14422.385
$ws.Range("N68").Value = -15920.385
$ws.Range("H71").Value = 14422.385
$ws.Range("J71").Value = 14422.385
$ws.Range("L71").Value = 43267.155
$ws.Range("N71").Value = -50755.155
$ws.Range("H74").Value = 23329.834
$ws.Range("J74").Value = 23329.834
$ws.Range("L74").Value = 23329.834
$ws.Range("N74").Value = -25077.834
$ws.Range("H77").Value = 23329.834
$ws.Range("J77").Value = 23329.834
$ws.Range("L77").Value = 69989.50199999999
$ws.Range("N77").Value = -78725.50199999999
$ws.Range("H86").Value = 2722.6667
$ws.Range("J86").Value = 2785.2856
$ws.Range("L86").Value = 2785.2856
$ws.Range("N86").Value = -5031.2856
$ws.Range("H89").Value = 2722.6667
$ws.Range("J89").Value = 2785.2856
$ws.Range("L89").Value = 13926.428
$ws.Range("N89").Value = -25158.428
$ws.Range("H99").Value = 8491.235000000001
$ws.Range("I99").Value = 2976
$ws.Range("J99").Value = 10789.25
$ws.Range("K99").Value = 2976
$ws.Range("L99").Value = 10789.25
$ws.Range("M99").Value = -1478
$ws.Range("N99").Value = -13785.25
$ws.Range("H126").Value = 8491.235000000001
$ws.Range("I126").Value = 2976
$ws.Range("J126").Value = 10789.25
$ws.Range("K126").Value = 8928
$ws.Range("L126").Value = 32367.75
$ws.Range("M126").Value = -6458
$ws.Range("N126").Value = -37307.75
$ws.Range("H132").Value = 3438.9355
$ws.Range("I132").Value = 3769.85
$ws.Range("J132").Value = 2837.2727
$ws.Range("K132").Value = 11309.55
$ws.Range("L132").Value = 8511.8181
$ws.Range("M132").Value = -8779.549999999999
$ws.Range("N132").Value = -13571.8181

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H101").Value = 3478.75
$ws.Range("J101").Value = 3478.75
$ws.Range("L101").Value = 10436.25
$ws.Range("N101").Value = -15304.25
$ws.Range("H109").Value = 2718.5557
$ws.Range("I109").Value = 1376.75
$ws.Range("J109").Value = 3792
$ws.Range("K109").Value = 4130.25
$ws.Range("L109").Value = 11376
$ws.Range("M109").Value = -3090.25
$ws.Range("N109").Value = -13456
$ws.Range("H132").Value = 4600.8
$ws.Range("I132").Value = 5001.3335
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 45012.0015
$ws.Range("L132").Value = 36000
$ws.Range("M132").Value = -42482.0015
$ws.Range("N132").Value = -41060

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = ""
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = ""
$ws.Range("H109").Value = 9285
$ws.Range("J109").Value = 9285
$ws.Range("L109").Value = 9285
$ws.Range("N109").Value = -11365

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2409.238
$ws.Range("I7").Value = 1858.25
$ws.Range("K7").Value = 1858.25
$ws.Range("M7").Value = -1746.25
$ws.Range("H22").Value = 700
$ws.Range("I22").Value = 700
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 700
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -405
$ws.Range("N22").Value = ""
$ws.Range("H27").Value = 700
$ws.Range("I27").Value = 700
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 700
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -593
$ws.Range("N27").Value = ""
$ws.Range("H40").Value = 1954.4
$ws.Range("I40").Value = 1438.9
$ws.Range("J40").Value = 2469.9
$ws.Range("K40").Value = 1438.9
$ws.Range("L40").Value = 2469.9
$ws.Range("M40").Value = -1302.9
$ws.Range("N40").Value = -2741.9
$ws.Range("H54").Value = 14815.167
$ws.Range("J54").Value = 14815.167
$ws.Range("L54").Value = 14815.167
$ws.Range("N54").Value = -16103.167
$ws.Range("H126").Value = 2409.238
$ws.Range("I126").Value = 1858.25
$ws.Range("K126").Value = 5574.75
$ws.Range("M126").Value = -3104.75
$ws.Range("H132").Value = 4235.316
$ws.Range("I132").Value = 4308.2
$ws.Range("K132").Value = 12924.6
$ws.Range("M132").Value = -10394.6
